$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'" + '29.391.64'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'" + '  -0.02%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'" + '1.842.70'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'" + '  -0.32%  '
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'" + '0.9983'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'" + '  -0.29%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'" + '240.30'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'" + '  -0.23%  '
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'" + '  +0.22%  '
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'" + '  -0.33%  '
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'" + '  -0.33%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'" + '0.2908'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'" + '  +0.47%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'" + '25.07'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'" + '  +2.72%  '
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'" + '0.07740'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'" + '  -0.11%  '
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'" + '1.843.96'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'" + '  -0.16%  '
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'" + '4.987'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'" + '  -0.52%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'" + '0.6792'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'" + '  -0.18%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'" + '0.00001021'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'" + '  -1.25%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'" + '82.11'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'" + '  -0.97%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'" + '6.282'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'" + '  +2.77%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'" + '29.378.02'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'" + '  -0.14%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'" + '229.41'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'" + '  +0.03%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'" + '12.35'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'" + '  +0.31%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'" + '0.9994'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'" + '  -0.25%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'" + '7.416'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'" + '  -0.25%  '
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'" + '  -0.22%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'" + '158.65'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'" + '  -0.33%  '
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'" + '  +1.06%  '
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'" + '0.1352'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'" + '  -2.41%  '
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'" + '  -0.95%  '
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'" + '0.06583'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'" + '  +15.71%  '
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'" + '1.437'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'" + '  +2.63%  '
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'" + '1.488'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'" + '  +0.79%  '
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'" + '4.072'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'" + '  -1.18%  '
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'" + '4.057'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'" + '  +0.15%  '
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'" + '  +0.97%  '
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "'" + '  -0.79%  '
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'" + '0.6985'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'" + '  +0.59%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'" + '2.577'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'" + '  -0.50%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'" + '0.01854'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'" + '  +1.92%  '
$c.Style = "Normal"

$c = $ws.Range("B38")
$c.Value = "'" + 'Maker'
$c.Style = "Normal"
$c = $ws.Range("C38")
$c.Value = "'" + 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'" + '1.250.05'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'" + '  -0.06%  '
$c.Style = "Normal"

$c = $ws.Range("B39")
$c.Value = "'" + 'MXToken'
$c.Style = "Normal"
$c = $ws.Range("C39")
$c.Value = "'" + 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'" + '2.817'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'" + '  -0.73%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'" + '6.802'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'" + '  +4.65%  '
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'" + '0.9325'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'" + '  +2.86%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'" + '0.9991'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'" + '  -0.24%  '
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'" + '1.995.49'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'" + '  -0.65%  '
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'" + '100.96'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'" + '  -0.41%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'" + '65.54'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'" + '  -0.45%  '
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'" + '  +2.68%  '
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'" + '  -0.30%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'" + '1.719'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'" + '  +3.70%  '
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'" + '9.042'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'" + '  +1.04%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'" + '0.1147'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'" + '  -1.03%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'" + '0.3904'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'" + '  -0.92%  '
$c.Style = "Normal"
